# Auto update Excel log
# Appends new sensor/alert log rows to four worksheets, matching the
# SeniorConnect sensor logging pipeline's latest readings.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$text)
    # Force the cell to be stored as Text so values such as "2026-02-01"
    # or "14:51:23" are not auto-converted into Excel date/time serials.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    # Drop back to the workbook's default (un-styled) cell style now that
    # the literal text value has been locked in, so no stray number format
    # is left applied to the cell.
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# ALERTS sheet: one new CRITICAL fall-detected alert (row 19)
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")

Set-TextCell $wsAlerts "A19" "2026-02-01"
Set-TextCell $wsAlerts "B19" "14:51:25"
Set-TextCell $wsAlerts "C19" "14:00"
$wsAlerts.Range("D19").Value = "Living Room"
$wsAlerts.Range("E19").Value = "CRITICAL"
$wsAlerts.Range("F19").Value = "FALL_DETECTED"

# ---------------------------------------------------------------------
# mmWave(InBed) sheet: fourteen new "In Bed" occupancy rows (22-35)
# ---------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

$inBedTimes = @(
    "14:51:23", "14:51:27", "14:51:29", "14:51:30", "14:51:31",
    "14:51:32", "14:51:33", "14:51:34", "14:51:35", "14:51:36",
    "14:51:37", "14:51:38", "14:51:39", "14:51:43"
)

for ($i = 0; $i -lt $inBedTimes.Count; $i++) {
    $r = 22 + $i
    Set-TextCell $wsInBed "A$r" "2026-02-01"
    Set-TextCell $wsInBed "B$r" $inBedTimes[$i]
    Set-TextCell $wsInBed "C$r" "14:00"
    $wsInBed.Range("D$r").Value = "Bedroom"
    $wsInBed.Range("E$r").Value = "In Bed"
    $wsInBed.Range("F$r").Value = "Occupied"
}

# ---------------------------------------------------------------------
# mmWave(BR) sheet: fourteen new breath-rate readings (22-35)
# ---------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

$brRows = @(
    @{ Time = "14:51:23"; Value = 1 },
    @{ Time = "14:51:28"; Value = 2 },
    @{ Time = "14:51:29"; Value = 9 },
    @{ Time = "14:51:31"; Value = 37 },
    @{ Time = "14:51:32"; Value = 3 },
    @{ Time = "14:51:33"; Value = 30 },
    @{ Time = "14:51:34"; Value = 3 },
    @{ Time = "14:51:35"; Value = 2 },
    @{ Time = "14:51:36"; Value = 13 },
    @{ Time = "14:51:37"; Value = 2 },
    @{ Time = "14:51:38"; Value = 37 },
    @{ Time = "14:51:39"; Value = 11 },
    @{ Time = "14:51:40"; Value = 2 },
    @{ Time = "14:51:43"; Value = 13 }
)

for ($i = 0; $i -lt $brRows.Count; $i++) {
    $r = 22 + $i
    $row = $brRows[$i]
    Set-TextCell $wsBR "A$r" "2026-02-01"
    Set-TextCell $wsBR "B$r" $row.Time
    Set-TextCell $wsBR "C$r" "14:00"
    $wsBR.Range("D$r").Value = "Bedroom"
    $wsBR.Range("E$r").Value = $row.Value
    $wsBR.Range("F$r").Value = "Occupied"
}

# ---------------------------------------------------------------------
# mmWave(HR) sheet: fourteen new heart-rate readings (22-35)
# ---------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

$hrRows = @(
    @{ Time = "14:51:23"; Value = 49 },
    @{ Time = "14:51:27"; Value = 50 },
    @{ Time = "14:51:29"; Value = 57 },
    @{ Time = "14:51:30"; Value = 85 },
    @{ Time = "14:51:31"; Value = 51 },
    @{ Time = "14:51:32"; Value = 78 },
    @{ Time = "14:51:33"; Value = 51 },
    @{ Time = "14:51:34"; Value = 50 },
    @{ Time = "14:51:35"; Value = 61 },
    @{ Time = "14:51:36"; Value = 50 },
    @{ Time = "14:51:37"; Value = 85 },
    @{ Time = "14:51:38"; Value = 59 },
    @{ Time = "14:51:39"; Value = 50 },
    @{ Time = "14:51:43"; Value = 61 }
)

for ($i = 0; $i -lt $hrRows.Count; $i++) {
    $r = 22 + $i
    $row = $hrRows[$i]
    Set-TextCell $wsHR "A$r" "2026-02-01"
    Set-TextCell $wsHR "B$r" $row.Time
    Set-TextCell $wsHR "C$r" "14:00"
    $wsHR.Range("D$r").Value = "Bedroom"
    $wsHR.Range("E$r").Value = $row.Value
    $wsHR.Range("F$r").Value = "Occupied"
}
